# Append 5 new rows (62-66) to the TestResults sheet, mirroring the pattern
# used by each existing 5-row block of test results, with a fresh
# "Bad Request" response timestamp for the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$postVideoRequest  = '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}'
$postVideoResponse = '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}'
$getVideoRequest   = 'GET https://www.videogamedb.uk:443/api/videogame/1'
$getVideoResponse  = '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}'
$putVideoResponse  = '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}'
$deleteVideoRequest  = 'DELETE https://www.videogamedb.uk:443/api/videogame/1'
$deleteVideoResponse = 'Video game deleted'
$invalidPostRequest  = '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}'
$invalidPostResponse = '{"timestamp":"2025-01-02T12:28:15.332+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}'

$rows = @(
    @("testPostVideo",      "Passed", $postVideoRequest,   $postVideoResponse),
    @("testGetVideoById",   "Passed", $getVideoRequest,     $getVideoResponse),
    @("testPutVideo",       "Passed", $postVideoRequest,    $putVideoResponse),
    @("testDeleteVideo",    "Passed", $deleteVideoRequest,  $deleteVideoResponse),
    @("testInvalidPostVideo","Passed", $invalidPostRequest, $invalidPostResponse)
)

$startRow = 62
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
